$wb = $excel.ActiveWorkbook

# ----- "Gen slack" sheet: update a couple of values -----
$wsGen = $wb.Worksheets.Item("Gen slack")
$wsGen.Range("A3").Value = 1
$wsGen.Range("A4").Value = 2
$wsGen.Range("F8").Select()

# ----- "Lines" sheet: renumber several bus references -----
$wsLines = $wb.Worksheets.Item("Lines")
$wsLines.Range("C3").Value = 3

$wsLines.Range("B4").Value = 4
$wsLines.Range("C4").Value = 5

$wsLines.Range("B5").Value = 5
$wsLines.Range("C5").Value = 6

$wsLines.Range("B6").Value = 6
$wsLines.Range("C6").Value = 7

$wsLines.Range("B7").Value = 7
$wsLines.Range("C7").Value = 3

$wsLines.Range("B8").Value = 8
$wsLines.Range("C8").Value = 9

$wsLines.Range("B9").Value = 9
$wsLines.Range("C9").Value = 10

# ----- "Load" sheet: move the selection (was the active tab) -----
$wsLoad = $wb.Worksheets.Item("Load")
$wsLoad.Range("F10").Select()

# ----- Make "Lines" the active tab/selection last, matching the saved view -----
$wsLines.Range("B10").Select()
